$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.124.73"
$ws.Range("E2").Value = "  -1.40%  "

# Row 3
$ws.Range("D3").Value = "2.244.78"
$ws.Range("E3").Value = "  -1.59%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.93"
$ws.Range("E5").Value = "  -1.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.51"
$ws.Range("E6").Value = "  -6.64%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.575"
$ws.Range("E7").Value = "  -3.47%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -7.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.29"
$ws.Range("E10").Value = "  -6.66%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0825"
$ws.Range("E11").Value = "  -2.51%  "

# Row 12
$ws.Range("E12").Value = "  -6.93%  "

# Row 13
$ws.Range("E13").Value = "  -3.03%  "

# Row 14
$ws.Range("D14").Value = "2.586.54"
$ws.Range("E14").Value = "  -1.70%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.845"
$ws.Range("E15").Value = "  -4.88%  "

# Row 16
$ws.Range("D16").Value = "2.246.28"
$ws.Range("E16").Value = "  -1.64%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.01"
$ws.Range("E17").Value = "  -4.56%  "

# Row 18
$ws.Range("D18").Value = "43.950.42"
$ws.Range("E18").Value = "  -1.56%  "

# Row 19
$ws.Range("E19").Value = "  -6.62%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("E20").Value = "  -2.83%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.33"
$ws.Range("E21").Value = "  -3.80%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.56"
$ws.Range("E22").Value = "  -1.72%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.37"
$ws.Range("E23").Value = "  -0.92%  "

# Row 24
$ws.Range("E24").Value = "  -7.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.03"
$ws.Range("E25").Value = "  -8.61%  "

# Row 26
$ws.Range("E26").Value = "  +0.46%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  -1.03%  "

# Row 28
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.06"
$ws.Range("E28").Value = "  -4.90%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").Value = "  -4.82%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.00"
$ws.Range("E30").Value = "  -8.43%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.06"
$ws.Range("E31").Value = "  -3.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.77"
$ws.Range("E32").Value = "  -4.39%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0839"
$ws.Range("E33").Value = "  -5.91%  "

# Row 34
$ws.Range("E34").Value = "  +4.72%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.67"
$ws.Range("E35").Value = "  -4.33%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"
$ws.Range("E36").Value = "  -7.06%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  -7.26%  "

# Row 38
$ws.Range("E38").Value = "  -3.06%  "

# Row 39
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.32"
$ws.Range("E39").Value = "  -2.44%  "

# Row 40
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.54"
$ws.Range("E40").Value = "  -10.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.99"
$ws.Range("E41").Value = "  -11.16%  "

# Row 42
$ws.Range("E42").Value = "  -6.26%  "

# Row 43
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("D44").Value = "1.702.05"
$ws.Range("E44").Value = "  -4.24%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "83.61"
$ws.Range("E45").Value = "  -4.62%  "

# Row 46
$ws.Range("E46").Value = "  -6.92%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.20"
$ws.Range("E47").Value = "  -5.45%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.94"
$ws.Range("E48").Value = "  -2.87%  "

# Row 49
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "56.69"
$ws.Range("E49").Value = "  -6.99%  "

# Row 50
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.14"
$ws.Range("E50").Value = "  -7.20%  "

# Row 51
$ws.Range("E51").Value = "  -6.71%  "
